$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.174.56"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").Value = "1.580.61"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "212.52"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  +6.21%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "26.14"
$ws.Range("E8").Value = "  +9.68%  "
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "1.806.70"
$ws.Range("D13").Value = "1.581.86"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "29.213.11"
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "3.71"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "62.76"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "236.71"
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "0.0₃0688"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("D25").Value = "154.15"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").Value = "15.14"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "0.0470"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "1.424.22"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").Value = "2.75"
$ws.Range("E37").Value = "  +6.15%  "
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").Value = "0.531"
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").Value = "52.97"
$ws.Range("E42").Value = "  +24.53%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "0.791"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").Value = "0.0472"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "64.54"
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "1.718.39"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").Value = "0.843"
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("D50").Value = "85.27"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -0.71%  "
